$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 275; this pushes the existing rows 275-293 down to 276-294,
# matching the rest of the diff (each old row's data shifted down by one row).
$ws.Rows.Item(275).Insert()

# Populate the newly inserted row 275 with the new record's data.
$ws.Range("A275").Value = 4
$ws.Range("B275").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C275").Value = "Los Lagos"
$ws.Range("D275").Value = 45013
$ws.Range("E275").Value = 10
$ws.Range("F275").Value = 100112009
$ws.Range("G275").Value = "Acelga"
$ws.Range("H275").Value = "Sin especificar"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 40
$ws.Range("K275").Value = 9000
$ws.Range("L275").Value = 9000
$ws.Range("M275").Value = 9000
$ws.Range("N275").Value = "$/docena de atados (12 kilos)"
$ws.Range("O275").Value = "Región de La Araucanía"
$ws.Range("P275").Value = 750
$ws.Range("Q275").Value = 12
$ws.Range("R275").Value = "Hortaliza"
